$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking stat columns (G:K) stay text for every data row,
# same as how they were already stored as text in the original workbook.
$ws.Range("A1:K6").NumberFormat = "@"

# -- Header row --
$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

# -- Row 2: vs Mumbai Indians, Abu Dhabi, Oct 28 2020 --
$ws.Range("A2").Value = " Abu Dhabi"
$ws.Range("B2").Value = " October 28 2020"
$ws.Range("C2").Value = "Mumbai won by 5 wickets (with 5 balls remaining)"
$ws.Range("D2").Value = "Royal Challengers Bangalore"
$ws.Range("E2").Value = "Mumbai Indians"
$ws.Range("F2").Value = "Gurkeerat Singh Mann "
$ws.Range("G2").Value = "14"
$ws.Range("H2").Value = "11"
$ws.Range("I2").Value = "2"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "127.27"

# -- Row 3: vs Rajasthan Royals, Dubai (DSC), Oct 17 2020 --
$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " October 17 2020"
$ws.Range("C3").Value = "RCB won by 7 wickets (with 2 balls remaining)"
$ws.Range("D3").Value = "Royal Challengers Bangalore"
$ws.Range("E3").Value = "Rajasthan Royals"
$ws.Range("F3").Value = "Gurkeerat Singh Mann "
$ws.Range("G3").Value = "19"
$ws.Range("H3").Value = "17"
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "111.76"

# -- Row 4: vs Sunrisers Hyderabad, Sharjah, Oct 31 2020 --
$ws.Range("A4").Value = " Sharjah"
$ws.Range("B4").Value = " October 31 2020"
$ws.Range("C4").Value = "Sunrisers won by 5 wickets (with 35 balls remaining)"
$ws.Range("D4").Value = "Royal Challengers Bangalore"
$ws.Range("E4").Value = "Sunrisers Hyderabad"
$ws.Range("F4").Value = "Gurkeerat Singh Mann "
$ws.Range("G4").Value = "15"
$ws.Range("H4").Value = "24"
$ws.Range("I4").Value = "1"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "62.50"

# -- Row 5: vs Kolkata Knight Riders, Abu Dhabi, Oct 21 2020 --
$ws.Range("A5").Value = " Abu Dhabi"
$ws.Range("B5").Value = " October 21 2020"
$ws.Range("C5").Value = "RCB won by 8 wickets (with 39 balls remaining)"
$ws.Range("D5").Value = "Royal Challengers Bangalore"
$ws.Range("E5").Value = "Kolkata Knight Riders"
$ws.Range("F5").Value = "Gurkeerat Singh Mann "
$ws.Range("G5").Value = "21"
$ws.Range("H5").Value = "26"
$ws.Range("I5").Value = "4"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "80.76"

# -- Row 6: vs Chennai Super Kings, Dubai (DSC), Oct 25 2020 --
$ws.Range("A6").Value = " Dubai (DSC)"
$ws.Range("B6").Value = " October 25 2020"
$ws.Range("C6").Value = "Super Kings won by 8 wickets (with 8 balls remaining)"
$ws.Range("D6").Value = "Royal Challengers Bangalore"
$ws.Range("E6").Value = "Chennai Super Kings"
$ws.Range("F6").Value = "Gurkeerat Singh Mann "
$ws.Range("G6").Value = "2"
$ws.Range("H6").Value = "2"
$ws.Range("I6").Value = "0"
$ws.Range("J6").Value = "0"
$ws.Range("K6").Value = "100.00"
